# Semana 43 de 2024 - update Poisson expected/observed/p-value table and
# insert a new "Sindrome de rubeola congenita" (evento 720) row before
# "Sarampion".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (shifts old rows 32-38 down to 33-39)
$ws.Rows(32).Insert()

# Row 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0.06

# Row 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 0.08

# Row 6
$ws.Range("D6").Value = 36

# Row 7
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.05

# Row 11
$ws.Range("C11").Value = 38
$ws.Range("D11").Value = 34
$ws.Range("E11").Value = 0.05

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.37

# Row 14
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 5

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0

# Row 18
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0.14

# Row 21
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.03

# Row 22
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0

# Row 24
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0.37

# Row 25
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0

# Row 26
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0

# Row 27
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 0.13

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0.18

# Row 31
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0.37

# Row 32
$ws.Range("A32").Value = "'720"
$ws.Range("A32").Style = "Normal"
$ws.Range("B32").Value = "Sindrome de rubeola congenita"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0

# Row 35
$ws.Range("D35").Value = 1

# Row 36
$ws.Range("C36").Value = 7
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0

# Row 37
$ws.Range("C37").Value = 10
$ws.Range("D37").Value = 3

# Row 38
$ws.Range("C38").Value = 8
$ws.Range("D38").Value = 2
$ws.Range("E38").Value = 0.01
